$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before E, shifting old E..P to F..Q
$ws.Columns("E").Insert()

# Header for new column E
$ws.Range("E1").Value = "pt_max"

# Values for new column E (rows 2-12)
$ws.Range("E2:E12").Value = 50

# Update selection to match the edit focus
$ws.Range("E2:E12").Select() | Out-Null
